$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D and E) for "ownTeam" and "oppTeam" before the
# existing "batsman" column, shifting the later columns to the right.
$ws.Columns("D:E").Insert()

# Insert a new row (row 2) for the Sept 24 2020 match, pushing the
# existing Sept 21 2020 match data down to row 3.
$ws.Rows("2:2").Insert()

# ---- Header row ----
$ws.Cells.Item(1,4).Value = "ownTeam"
$ws.Cells.Item(1,5).Value = "oppTeam"

# ---- Row 2: new match data (Sept 24 2020 vs Kings XI Punjab) ----
$ws.Cells.Item(2,1).Value = " Dubai (DSC)"
$ws.Cells.Item(2,2).Value = " September 24 2020"
$ws.Cells.Item(2,3).Value = "Kings XI won by 97 runs"
$ws.Cells.Item(2,4).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(2,5).Value = "Kings XI Punjab"
$ws.Cells.Item(2,6).Value = "Josh Philippe" + [char]0x00A0 + "†"
$ws.Cells.Item(2,7).Value = "'0"
$ws.Cells.Item(2,8).Value = "'3"
$ws.Cells.Item(2,9).Value = "'0"
$ws.Cells.Item(2,10).Value = "'0"
$ws.Cells.Item(2,11).Value = "'0.00"

# ---- Row 3: existing match data (Sept 21 2020 vs Sunrisers Hyderabad) ----
$ws.Cells.Item(3,4).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(3,5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(3,7).Value = "'1"
$ws.Cells.Item(3,8).Value = "'2"
$ws.Cells.Item(3,9).Value = "'0"
$ws.Cells.Item(3,10).Value = "'0"
$ws.Cells.Item(3,11).Value = "'50.00"
